$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dog bone -> add Amazon.de link in E26
$ws.Hyperlinks.Add($ws.Range("E26"), "https://www.amazon.de/-/en/dp/B08FMJXFCH?psc=1&ref=ppx_yo2ov_dt_b_product_details", "", "", "https://www.amazon.de/-/en/dp/B08FMJXFCH?psc=1&ref=ppx_yo2ov_dt_b_product_details")

# Pinion -> add Absima shop link in E29
$ws.Hyperlinks.Add($ws.Range("E29"), "https://www.absima.shop/pp/alu-pinion-32dp/module0-8-20T.htm?shop=absima_en&SessionId=&a=article&ProdNr=2310348&t=19114&c=19132&p=19132", "", "", "https://www.absima.shop/pp/alu-pinion-32dp/module0-8-20T.htm?shop=absima_en&SessionId=&a=article&ProdNr=2310348&t=19114&c=19132&p=19132")

# Battery -> add Amazon.de link in E30
$ws.Hyperlinks.Add($ws.Range("E30"), "https://www.amazon.de/-/en/dp/B08X4GF9DK?psc=1&ref=ppx_yo2ov_dt_b_product_details", "", "", "https://www.amazon.de/-/en/dp/B08X4GF9DK?psc=1&ref=ppx_yo2ov_dt_b_product_details")

# Match the style used by the other "Link" column cells (E2 etc.) rather than
# Excel's default built-in Hyperlink style that .Hyperlinks.Add() applies.
$ws.Range("E2").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E30").PasteSpecial(-4122)
$excel.CutCopyMode = 0
